# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Leve market-profit tables
# across all eight crafting-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 277.86667
$ws.Range("I11").Value = 277.86667
$ws.Range("K11").Value = 277.86667
$ws.Range("M11").Value = -137.86667
$ws.Range("H58").Value = 916.2857
$ws.Range("I58").Value = 569
$ws.Range("K58").Value = 1707
$ws.Range("M58").Value = -1557
$ws.Range("H131").Value = 366.66666
$ws.Range("I131").Value = 366.66666
$ws.Range("K131").Value = 1099.99998
$ws.Range("M131").Value = 3940.00002
$ws.Range("H138").Value = 1217.6428
$ws.Range("H141").Value = 2492.6667
$ws.Range("I141").Value = 1810.1818
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 5430.5454
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -250.5454
$ws.Range("N141").Value = -40360

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 2336
$ws.Range("I36").Value = 1803.2
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 1803.2
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = -1457.2
$ws.Range("N36").Value = -5692
$ws.Range("H74").Value = 6712
$ws.Range("I74").Value = 6712
$ws.Range("K74").Value = 6712
$ws.Range("M74").Value = -5838
$ws.Range("H77").Value = 6712
$ws.Range("I77").Value = 6712
$ws.Range("K77").Value = 33560
$ws.Range("M77").Value = -29192
$ws.Range("H122").Value = 1874.5
$ws.Range("I122").Value = 1874.5
$ws.Range("K122").Value = 5623.5
$ws.Range("M122").Value = -3173.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3749.75
$ws.Range("I20").Value = 1999.5
$ws.Range("J20").Value = 5500
$ws.Range("K20").Value = 1999.5
$ws.Range("L20").Value = 5500
$ws.Range("M20").Value = -1752.5
$ws.Range("N20").Value = -5994
$ws.Range("H22").Value = 421.42856
$ws.Range("I22").Value = 394.11765
$ws.Range("J22").Value = 537.5
$ws.Range("K22").Value = 394.11765
$ws.Range("L22").Value = 537.5
$ws.Range("M22").Value = -221.11765
$ws.Range("N22").Value = -883.5
$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 1000
$ws.Range("M99").Value = 498

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1669
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = $null
$ws.Range("H34").Value = 1669
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = $null
$ws.Range("H58").Value = 4869.1113
$ws.Range("I58").Value = 6334.6665
$ws.Range("J58").Value = 1938
$ws.Range("K58").Value = 6334.6665
$ws.Range("L58").Value = 1938
$ws.Range("M58").Value = -6131.6665
$ws.Range("N58").Value = -2344
$ws.Range("H99").Value = 4361.125
$ws.Range("I99").Value = 3974.75
$ws.Range("J99").Value = 4747.5
$ws.Range("K99").Value = 3974.75
$ws.Range("L99").Value = 4747.5
$ws.Range("M99").Value = -2476.75
$ws.Range("N99").Value = -7743.5
$ws.Range("H122").Value = 1675
$ws.Range("I122").Value = 1400
$ws.Range("K122").Value = 4200
$ws.Range("M122").Value = -1750
$ws.Range("H126").Value = 4361.125
$ws.Range("I126").Value = 3974.75
$ws.Range("J126").Value = 4747.5
$ws.Range("K126").Value = 11924.25
$ws.Range("L126").Value = 14242.5
$ws.Range("M126").Value = -9454.25
$ws.Range("N126").Value = -19182.5
$ws.Range("H136").Value = 4869.1113
$ws.Range("I136").Value = 6334.6665
$ws.Range("J136").Value = 1938
$ws.Range("K136").Value = 19003.9995
$ws.Range("L136").Value = 5814
$ws.Range("M136").Value = -16453.9995
$ws.Range("N136").Value = -10914

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 196.3
$ws.Range("J12").Value = 257.57144
$ws.Range("L12").Value = 772.71432
$ws.Range("N12").Value = -1118.71432
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").Value = $null
$ws.Range("H117").Value = 640.125
$ws.Range("I117").Value = 543
$ws.Range("J117").Value = 931.5
$ws.Range("K117").Value = 1629
$ws.Range("L117").Value = 2794.5
$ws.Range("M117").Value = 1813
$ws.Range("N117").Value = -9678.5
$ws.Range("H122").Value = 738.4
$ws.Range("I122").Value = 673
$ws.Range("K122").Value = 6057
$ws.Range("M122").Value = -3607

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("L21").Value = $null
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("L30").Value = $null
$ws.Range("H31").Value = 1274.75
$ws.Range("I31").Value = 1274.75
$ws.Range("K31").Value = 1274.75
$ws.Range("M31").Value = -982.75
$ws.Range("H37").Value = 1274.75
$ws.Range("I37").Value = 1274.75
$ws.Range("K37").Value = 1274.75
$ws.Range("M37").Value = -997.75
$ws.Range("H44").Value = 20333
$ws.Range("J44").Value = 20333
$ws.Range("L44").Value = 20333
$ws.Range("N44").Value = -21525
$ws.Range("H70").Value = 7999
$ws.Range("I70").Value = 4000
$ws.Range("K70").Value = 4000
$ws.Range("M70").Value = -3730
$ws.Range("H73").Value = 7999
$ws.Range("I73").Value = 4000
$ws.Range("K73").Value = 4000
$ws.Range("M73").Value = -3064
$ws.Range("H101").Value = 36000
$ws.Range("J101").Value = 36000
$ws.Range("L101").Value = 36000
$ws.Range("N101").Value = -42490
$ws.Range("H102").Value = 1396.4
$ws.Range("I102").Value = 1396.4
$ws.Range("K102").Value = 1396.4
$ws.Range("M102").Value = 225.5999999999999
$ws.Range("H122").Value = 10420020
$ws.Range("I122").Value = 13891515
$ws.Range("K122").Value = 41674545
$ws.Range("M122").Value = -41672095

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1103.6666
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = $null
$ws.Range("H46").Value = 2200.611
$ws.Range("I46").Value = 2245
$ws.Range("J46").Value = 2130.8572
$ws.Range("K46").Value = 2245
$ws.Range("L46").Value = 2130.8572
$ws.Range("M46").Value = -2057
$ws.Range("N46").Value = -2506.8572
$ws.Range("H61").Value = 7434.5
$ws.Range("I61").Value = 7434.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 7434.5
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = -7232.5
$ws.Range("M61").Value = $null
$ws.Range("H82").Value = 1622.8334
$ws.Range("J82").Value = 1800.5
$ws.Range("L82").Value = 1800.5
$ws.Range("N82").Value = -2522.5
$ws.Range("H85").Value = 1622.8334
$ws.Range("J85").Value = 1800.5
$ws.Range("L85").Value = 1800.5
$ws.Range("N85").Value = -4296.5
$ws.Range("H113").Value = 7434.5
$ws.Range("I113").Value = 7434.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 7434.5
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = -5264.5
$ws.Range("M113").Value = $null
$ws.Range("H122").Value = 3504
$ws.Range("I122").Value = 3504
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10512
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = -8062
$ws.Range("M122").Value = $null
$ws.Range("H136").Value = 2724.0715
$ws.Range("I136").Value = 2625.5386
$ws.Range("K136").Value = 7876.6158
$ws.Range("M136").Value = -5326.6158

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("N113").Value = 0
$ws.Range("L113").Value = $null
$ws.Range("M113").Value = $null
$ws.Range("H136").Value = 1330.5714
$ws.Range("I136").Value = 1386.3846
$ws.Range("J136").Value = 605
$ws.Range("K136").Value = 4159.1538
$ws.Range("L136").Value = 1815
$ws.Range("M136").Value = -1609.1538
$ws.Range("N136").Value = -6915
